$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 26 (shifts rows 26-69 down to 27-70)
$ws.Rows.Item(26).Insert()

# Populate the newly inserted row 26 with data
$ws.Cells.Item(26, 1).Value = 3
$ws.Cells.Item(26, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(26, 3).Value = "Coquimbo"
$ws.Cells.Item(26, 4).Value = 44477
$ws.Cells.Item(26, 5).Value = 5
$ws.Cells.Item(26, 6).Value = 100112026
$ws.Cells.Item(26, 7).Value = "Haba"
$ws.Cells.Item(26, 8).Value = "Sin especificar"
$ws.Cells.Item(26, 9).Value = "Primera"
$ws.Cells.Item(26, 10).Value = 85
$ws.Cells.Item(26, 11).Value = 9000
$ws.Cells.Item(26, 12).Value = 10000
$ws.Cells.Item(26, 13).Value = 9529
$ws.Cells.Item(26, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(26, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(26, 16).Value = 381
$ws.Cells.Item(26, 17).Value = 25
$ws.Cells.Item(26, 18).Value = "Hortaliza"
